$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New integer reward values replacing the old decimal ones (B2:D13)
$values = @(
    @(20, -10, -20),
    @(20, -10, -20),
    @(20, -10, -20),
    @(-20, -20, -20),
    @(20, 0, -20),
    @(20, 0, -20),
    @(20, 0, -20),
    @(10, 10, -20),
    @(20, 10, -20),
    @(20, 10, -20),
    @(20, 10, -20),
    @(20, 20, -20)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i][0]
    $ws.Cells.Item($row, 3).Value = $values[$i][1]
    $ws.Cells.Item($row, 4).Value = $values[$i][2]
}

# Update the active selection
$ws.Range("C31").Select()
